$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (E1) into the new
# header cell F1, then set its text to "Modelo".
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Add the model name value in the corresponding data row (F2), using
# default (unstyled) formatting like the other data cells.
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
